$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 300
$ws.Range("F9").Value = 300

$ws.Range("F4").Select()
